$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a simple bug/feature tracker:
#  - rows 9-18 (col A = description, col B = status) get marked "corregido"
#    for items that have since been fixed
#  - row 19's old item was dropped from this list
#  - rows 20-36 are the new numbered backlog (col A = item #, col B = description)

# Remove the old row 19 entry entirely (no longer present in the updated list)
$ws.Range("A19").ClearContents()

# NOTE: cells below are intentionally written in this particular order (not
# simple row order) so that newly-introduced text values land in the shared
# string table in the same sequence as the source workbook.
$ws.Range("B27").Value2 = "los numeros de camiseta tienen que poder modificarse antes de iniciar el juego"
$ws.Range("A17").Value2 = "falta un salir de la visualizacion del partido en espectador para volver a la pagina anterior o al menu o a los partidos (decidilo vos)"
$ws.Range("B20").Value2 = "falta un `"suspender`" partido para que el partido quede guardado todos los datos al momento de la suspension, al momento de suspenderse debe abrir un cuadro de observaciones explicando el por qué de la suspension"
$ws.Range("B21").Value2 = "al momento de agregar una falta, deberia abrirse un menu flotante que indique que tipo de falta es (personal, tecnica, antideportiva, descalificadora) y (salvo la falta tecnica) marcar la cantidad de lanzamientos que debe hacer 0,1,2,3."
$ws.Range("B23").Value2 = "al iniciar el juego, pero antes de elegir a los 5 titulares, los dos equipos deberian poder elegir a los 12 jugadores citados para el juego. Ya que puede haber equipos que tienen mas de 12 jugadores anotados para jugar el torneo"
$ws.Range("B26").Value2 = "tambien deben aparecer las correcciones (descuentos)"
$ws.Range("B28").Value2 = "el entrenador tambien debe aparecer en algun lugar (podria ser debajo de los 5 titulares, el cual tambien puede recibir falta tecnica o descalificatoria"
$ws.Range("A15").Value2 = "al apretar Partidos en Vivo, entra en todos los partidos, deberia entrar directamente a los partidos en vivo"
$ws.Range("A16").Value2 = "en el index, los partidos en vivo y los ultimos resultados no estan apareciendo"
$ws.Range("B29").Value2 = "para qué sirve el boton Cargar Partido si me lleva a la seccion Partidos? Al igual que el boton Partidos. Es para que despues use el administrador de torneos?"
$ws.Range("B30").Value2 = "falta un salir de la carga  del partido en el planiller para volver a la pagina anterior o al menu o a los partidos (decidilo vos)"
$ws.Range("B31").Value2 = "si el jugador está con cero faltas se le pueden seguir descontando y las sacas de las falta de equipo"
$ws.Range("B32").Value2 = "en el celular hay que arreglar la adaptabilidad y ver por qué en la app no gira la pantalla, ya que es mucho mas comodo usarlo de esta forma"

# Remaining cells (unchanged values, or values already present in the shared
# string table) can be written in plain top-to-bottom, left-to-right order.
$ws.Range("A2").Value2 = "San Lorenzo: "
$ws.Range("B2").Value2 = "109ff4e8-0889-4b03-8c5d-80c91d60fc06"
$ws.Range("A3").Value2 = "Obras: "
$ws.Range("B3").Value2 = "90e53bf1-b65d-427d-b332-82b436082cd9"
$ws.Range("A4").Value2 = "Ferro: "
$ws.Range("B4").Value2 = "a3f25df9-fc00-4bf5-817f-1afd4721d842"
$ws.Range("A5").Value2 = " Boca: "
$ws.Range("B5").Value2 = "aeefd149-3ce7-401a-9773-5ee480f6e81a"
$ws.Range("A9").Value2 = "hay un error al solicitar minutos en el 1er tiempo, ya descuenta el 3er minuto del 2do tiempo"
$ws.Range("B9").Value2 = "corregido"
$ws.Range("A10").Value2 = "se pueden descontar puntos a jugadores que no hicieron puntos"
$ws.Range("B10").Value2 = "corregido"
$ws.Range("A11").Value2 = "el partido puede terminar empatado"
$ws.Range("B11").Value2 = "corregido"
$ws.Range("A12").Value2 = "cuando salis del partido y volves a entrar, te saca los titulares y al no haber cambios disponibles, no podes volver a meterlos, podriamos poner que los casilleros esten vacios pero cliqueables para poder ingresar (si no se puede corregir el error)"
$ws.Range("B12").Value2 = "corregido"
$ws.Range("A13").Value2 = "el recuadro de cada partido debe tener un fondo un poco mas oscuro para queden mejor separados entre si"
$ws.Range("B13").Value2 = "corregido"
$ws.Range("A14").Value2 = "los partidos deberian estar ordenados desde el mas proximo al mas lejano y despues los ya jugados"
$ws.Range("B14").Value2 = "corregido"
$ws.Range("B15").Value2 = "corregido"
$ws.Range("B16").Value2 = "corregido"
$ws.Range("B17").Value2 = "corregido"
$ws.Range("A18").Value2 = "el jugador que hace 5 faltas esta bien que no pueda tener mas acciones, pero debe poder hacer sustitucion para que entre otro jugador suplente"
$ws.Range("B18").Value2 = "corregido"
$ws.Range("A20").Value2 = 1
$ws.Range("A21").Value2 = 2
$ws.Range("A22").Value2 = 3
$ws.Range("B22").Value2 = "en caso de sumar 2 faltas tecnicas, 2 faltas antideportivas o 1 falta tecn y 1 desc, el jugador queda descalificado de ese partido debe aparecer un GD"
$ws.Range("A23").Value2 = 4
$ws.Range("A24").Value2 = 5
$ws.Range("B24").Value2 = "al momento de cargar jugadores nuevos se tiene que poder marcar al mismo como jugador Refuerzo, el cual va a tener reglas especiales en cuanto a la cantidad de cuartos que puede jugar"
$ws.Range("A25").Value2 = 6
$ws.Range("B25").Value2 = "en los logacciones deberia aparecer que terminó el cuarto y al momento de hacer la falta debe informar, tipo de falta, si da lanzamientos y cantidad de faltas de ese jugador"
$ws.Range("A26").Value2 = 7
$ws.Range("A27").Value2 = 8
$ws.Range("A28").Value2 = 9
$ws.Range("A29").Value2 = 10
$ws.Range("A30").Value2 = 11
$ws.Range("A31").Value2 = 12
$ws.Range("A32").Value2 = 13
$ws.Range("A33").Value2 = 14
$ws.Range("A34").Value2 = 15
$ws.Range("A35").Value2 = 16
$ws.Range("A36").Value2 = 17

# Update view state to match target (scroll position + active selection)
$ws.Activate()
$excel.Goto($ws.Range("A15"), $true)
$ws.Range("B33").Select()
